$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Success" sheet's C column was mislabeled "location" but actually held
# city names, and was missing a separate state column. Insert a new column
# before the old D ("descriptor") column to hold state data, then relabel
# and repopulate the columns so the final layout is:
#   A: amout | B: year | C: city | D: state | E: descriptor | F: type
$ws.Columns.Item(4).Insert()

# Rename the old "location" header to "city"
$ws.Cells.Item(1,3).Value = "city"

# Populate the new "state" column (header + 18 data rows)
$ws.Cells.Item(1,4).Value = "state"
$ws.Cells.Item(2,4).Value = "Illinois"
$ws.Cells.Item(3,4).Value = "Illinois"
$ws.Cells.Item(4,4).Value = "Texas"
$ws.Cells.Item(5,4).Value = "Tennessee"
$ws.Cells.Item(6,4).Value = "Texas"
$ws.Cells.Item(7,4).Value = "California"
$ws.Cells.Item(8,4).Value = "Florida"
$ws.Cells.Item(9,4).Value = "North Carolina"
$ws.Cells.Item(10,4).Value = "New York"
$ws.Cells.Item(11,4).Value = "Colorado"
$ws.Cells.Item(12,4).Value = "New Jersey"
$ws.Cells.Item(13,4).Value = "New York"
$ws.Cells.Item(14,4).Value = "California"
$ws.Cells.Item(15,4).Value = "Texas"
$ws.Cells.Item(16,4).Value = "California"
$ws.Cells.Item(17,4).Value = "California"
$ws.Cells.Item(18,4).Value = "California"

# Match the author's manual column width for the new state column
$ws.Columns.Item(4).ColumnWidth = 14.14

# Leave the cursor where the author left it when finishing up
$null = $ws.Range("D19").Select()
